$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (tab) from "units" to "Sheet1"
$ws.Name = "Sheet1"

# Capitalize a few mis-cased character names (data overwrite)
$ws.Range("A3").Value = "Shamann"
$ws.Range("A4").Value = "Astar"
$ws.Range("A5").Value = "Teo"
$ws.Range("A11").Value = "Tushen"
